$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # G2: was numeric 29.9 -> becomes inline string "不可售"
    $ws.Range("G2").Value = "不可售"

    # F3: was 990 -> becomes 991
    $ws.Range("F3").Value = 991
}
